# Updated Hardware section to include Large Cable Transmission System,
# Added in Publication Experimental Data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "Chain (2ft)" item (row 6) to "05B Chain (2ft)"
$ws.Range("A6").Value = "05B Chain (2ft)"

# Append a new line item: Bowden Tube (Large Cable Transmission System part)
$ws.Range("A20").Value = "Bowden Tube"
$ws.Range("B20").Value = "https://www.amazon.com/Jagwire-Sport-Housing-Slick-Lube-Titanium/dp/B085NBZMJS/ref=sr_1_1?crid=1A5WX5ADQYM0Y&dib=eyJ2IjoiMSJ9.CWV7EelBoN67bHKqG_VMGaAHwKq3lWAqkaCYopLdT43GyH4CDNeyWoQ_bFV_YrbZmhmwmsofP69GRzCCWYW_ULIkinZgZrdky8EGo_FPRa2GDLIPcrjwSu8T1nDFsZ03wyuLLatilRsdmpFkqvecV8S7AOhy1XjGVW6Ztcl1kgVL4_2zjOWbooP9z_kj4elJxdrMd7yL-uhr5ZcYm6F_Z725qnNy_c0-wNfaDt_xxhvLSabJzp2Ta9HzGpozVxtWGMo89NRT0qD_9iY1NxHHBedJp86w9TXZ0OfMm9s45WQ.7sQYhAOJzDOGkHvrSoGgv6F0CIzum6_7bjb1Zjv8zzc&dib_tag=se&keywords=Jagwire+Brake+Housing+CGX-SL+Slick-Lube+5+mm+%2810+m%29&qid=1726249759&s=sporting-goods&sprefix=jagwire+brake+housing+cgx-sl+slick-lube+5+mm+10+m+%2Csporting%2C99&sr=1-1"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 61.97

# Update the selected cell to match the saved session state
[void]$ws.Range("E33").Select()
